$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("QA3-short")

$ws.Range("C2").Value = "John"
$ws.Range("C8").Value = "Turner"

$ws.Range("F19").Select()
